$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = "BTC"
$ws.Cells.Item(2, 3).Value = "Bitcoin"
$ws.Cells.Item(2, 4).Value = 30343
$ws.Cells.Item(2, 5).Value = 589544684001
$ws.Cells.Item(2, 6).Value = 6019402119
$ws.Cells.Item(2, 7).Value = 0.00332

$ws.Cells.Item(3, 2).Value = "ETH"
$ws.Cells.Item(3, 3).Value = "Ethereum"
$ws.Cells.Item(3, 4).Value = 1935.52
$ws.Cells.Item(3, 5).Value = 232618512297
$ws.Cells.Item(3, 6).Value = 4230895730
$ws.Cells.Item(3, 7).Value = -0.24192

$ws.Cells.Item(4, 2).Value = "USDT"
$ws.Cells.Item(4, 3).Value = "Tether"
$ws.Cells.Item(4, 4).Value = 1.001
$ws.Cells.Item(4, 5).Value = 83669568053
$ws.Cells.Item(4, 6).Value = 12470789930
$ws.Cells.Item(4, 7).Value = -0.01269

$ws.Cells.Item(5, 2).Value = "XRP"
$ws.Cells.Item(5, 3).Value = "XRP"
$ws.Cells.Item(5, 4).Value = 0.751406
$ws.Cells.Item(5, 5).Value = 39485209298
$ws.Cells.Item(5, 6).Value = 3054239962
$ws.Cells.Item(5, 7).Value = 3.62814

$ws.Cells.Item(6, 2).Value = "BNB"
$ws.Cells.Item(6, 3).Value = "BNB"
$ws.Cells.Item(6, 4).Value = 250.45
$ws.Cells.Item(6, 5).Value = 39030804145
$ws.Cells.Item(6, 6).Value = 418541878
$ws.Cells.Item(6, 7).Value = -0.24334

$ws.Cells.Item(7, 2).Value = "USDC"
$ws.Cells.Item(7, 3).Value = "USD Coin"
$ws.Cells.Item(7, 4).Value = 1.001
$ws.Cells.Item(7, 5).Value = 27299501274
$ws.Cells.Item(7, 6).Value = 1648574761
$ws.Cells.Item(7, 7).Value = 0.03536

$ws.Cells.Item(8, 2).Value = "STETH"
$ws.Cells.Item(8, 3).Value = "Lido Staked Ether"
$ws.Cells.Item(8, 4).Value = 1934.41
$ws.Cells.Item(8, 5).Value = 14974134987
$ws.Cells.Item(8, 6).Value = 3991830
$ws.Cells.Item(8, 7).Value = -0.22958

$ws.Cells.Item(9, 2).Value = "ADA"
$ws.Cells.Item(9, 3).Value = "Cardano"
$ws.Cells.Item(9, 4).Value = 0.323311
$ws.Cells.Item(9, 5).Value = 11344300187
$ws.Cells.Item(9, 6).Value = 233325738
$ws.Cells.Item(9, 7).Value = -2.9888

$ws.Cells.Item(10, 2).Value = "SOL"
$ws.Cells.Item(10, 3).Value = "Solana"
$ws.Cells.Item(10, 4).Value = 28.14
$ws.Cells.Item(10, 5).Value = 11341103614
$ws.Cells.Item(10, 6).Value = 841124593
$ws.Cells.Item(10, 7).Value = -2.53674

$ws.Cells.Item(11, 2).Value = "DOGE"
$ws.Cells.Item(11, 3).Value = "Dogecoin"
$ws.Cells.Item(11, 4).Value = 0.07169
$ws.Cells.Item(11, 5).Value = 10037623091
$ws.Cells.Item(11, 6).Value = 683926724
$ws.Cells.Item(11, 7).Value = -3.05223

$ws.Cells.Item(12, 2).Value = "MATIC"
$ws.Cells.Item(12, 3).Value = "Polygon"
$ws.Cells.Item(12, 4).Value = 0.792743
$ws.Cells.Item(12, 5).Value = 7391342244
$ws.Cells.Item(12, 6).Value = 195075352
$ws.Cells.Item(12, 7).Value = -2.69314

$ws.Cells.Item(13, 2).Value = "TRX"
$ws.Cells.Item(13, 3).Value = "TRON"
$ws.Cells.Item(13, 4).Value = 0.080164
$ws.Cells.Item(13, 5).Value = 7195444700
$ws.Cells.Item(13, 6).Value = 199620290
$ws.Cells.Item(13, 7).Value = -1.26956

$ws.Cells.Item(14, 2).Value = "LTC"
$ws.Cells.Item(14, 3).Value = "Litecoin"
$ws.Cells.Item(14, 4).Value = 94.73999999999999
$ws.Cells.Item(14, 5).Value = 6949430314
$ws.Cells.Item(14, 6).Value = 532406213
$ws.Cells.Item(14, 7).Value = -0.12233

$ws.Cells.Item(15, 2).Value = "DOT"
$ws.Cells.Item(15, 3).Value = "Polkadot"
$ws.Cells.Item(15, 4).Value = 5.4
$ws.Cells.Item(15, 5).Value = 6782096112
$ws.Cells.Item(15, 6).Value = 116415831
$ws.Cells.Item(15, 7).Value = -1.68405

$ws.Cells.Item(16, 2).Value = "AVAX"
$ws.Cells.Item(16, 3).Value = "Avalanche"
$ws.Cells.Item(16, 4).Value = 14.62
$ws.Cells.Item(16, 5).Value = 5049286693
$ws.Cells.Item(16, 6).Value = 186434379
$ws.Cells.Item(16, 7).Value = -2.65001

$ws.Cells.Item(17, 2).Value = "BCH"
$ws.Cells.Item(17, 3).Value = "Bitcoin Cash"
$ws.Cells.Item(17, 4).Value = 254.1
$ws.Cells.Item(17, 5).Value = 4945184316
$ws.Cells.Item(17, 6).Value = 340341379
$ws.Cells.Item(17, 7).Value = 0.55036

$ws.Cells.Item(18, 2).Value = "WBTC"
$ws.Cells.Item(18, 3).Value = "Wrapped Bitcoin"
$ws.Cells.Item(18, 4).Value = 30332
$ws.Cells.Item(18, 5).Value = 4798778662
$ws.Cells.Item(18, 6).Value = 50688406
$ws.Cells.Item(18, 7).Value = 0.007990000000000001

$ws.Cells.Item(19, 2).Value = "SHIB"
$ws.Cells.Item(19, 3).Value = "Shiba Inu"
$ws.Cells.Item(19, 4).Value = 0.00000807
$ws.Cells.Item(19, 5).Value = 4754482965
$ws.Cells.Item(19, 6).Value = 145855361
$ws.Cells.Item(19, 7).Value = -2.64612

$ws.Cells.Item(20, 2).Value = "UNI"
$ws.Cells.Item(20, 3).Value = "Uniswap"
$ws.Cells.Item(20, 4).Value = 5.79
$ws.Cells.Item(20, 5).Value = 4365347946
$ws.Cells.Item(20, 6).Value = 74778790
$ws.Cells.Item(20, 7).Value = -1.75419

$ws.Cells.Item(21, 2).Value = "DAI"
$ws.Cells.Item(21, 3).Value = "Dai"
$ws.Cells.Item(21, 4).Value = 1
$ws.Cells.Item(21, 5).Value = 4287102528
$ws.Cells.Item(21, 6).Value = 41140048
$ws.Cells.Item(21, 7).Value = 0.03762

$ws.Cells.Item(22, 2).Value = "BUSD"
$ws.Cells.Item(22, 3).Value = "Binance USD"
$ws.Cells.Item(22, 4).Value = 1.001
$ws.Cells.Item(22, 5).Value = 3972150668
$ws.Cells.Item(22, 6).Value = 1589487997
$ws.Cells.Item(22, 7).Value = 0.03933

$ws.Cells.Item(23, 2).Value = "LEO"
$ws.Cells.Item(23, 3).Value = "LEO Token"
$ws.Cells.Item(23, 4).Value = 4.01
$ws.Cells.Item(23, 5).Value = 3723007671
$ws.Cells.Item(23, 6).Value = 811065
$ws.Cells.Item(23, 7).Value = 0.84906

$ws.Cells.Item(24, 2).Value = "XLM"
$ws.Cells.Item(24, 3).Value = "Stellar"
$ws.Cells.Item(24, 4).Value = 0.134775
$ws.Cells.Item(24, 5).Value = 3664590854
$ws.Cells.Item(24, 6).Value = 200501742
$ws.Cells.Item(24, 7).Value = 2.27628

$ws.Cells.Item(25, 2).Value = "LINK"
$ws.Cells.Item(25, 3).Value = "Chainlink"
$ws.Cells.Item(25, 4).Value = 6.84
$ws.Cells.Item(25, 5).Value = 3543826185
$ws.Cells.Item(25, 6).Value = 141533515
$ws.Cells.Item(25, 7).Value = -1.92002

$ws.Cells.Item(26, 2).Value = "XMR"
$ws.Cells.Item(26, 3).Value = "Monero"
$ws.Cells.Item(26, 4).Value = 164.43
$ws.Cells.Item(26, 5).Value = 2986834546
$ws.Cells.Item(26, 6).Value = 57203738
$ws.Cells.Item(26, 7).Value = 0.9755

$ws.Cells.Item(27, 2).Value = "TUSD"
$ws.Cells.Item(27, 3).Value = "TrueUSD"
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(27, 5).Value = 2845494888
$ws.Cells.Item(27, 6).Value = 894107610
$ws.Cells.Item(27, 7).Value = 0.07815999999999999

$ws.Cells.Item(28, 2).Value = "ATOM"
$ws.Cells.Item(28, 3).Value = "Cosmos Hub"
$ws.Cells.Item(28, 4).Value = 9.6
$ws.Cells.Item(28, 5).Value = 2810842987
$ws.Cells.Item(28, 6).Value = 64113749
$ws.Cells.Item(28, 7).Value = -2.48926

$ws.Cells.Item(29, 2).Value = "ETC"
$ws.Cells.Item(29, 3).Value = "Ethereum Classic"
$ws.Cells.Item(29, 4).Value = 19.13
$ws.Cells.Item(29, 5).Value = 2718361584
$ws.Cells.Item(29, 6).Value = 87657147
$ws.Cells.Item(29, 7).Value = -1.53924

$ws.Cells.Item(30, 2).Value = "OKB"
$ws.Cells.Item(30, 3).Value = "OKB"
$ws.Cells.Item(30, 4).Value = 44.37
$ws.Cells.Item(30, 5).Value = 2663179561
$ws.Cells.Item(30, 6).Value = 2765240
$ws.Cells.Item(30, 7).Value = 1.96364

$ws.Cells.Item(31, 2).Value = "LDO"
$ws.Cells.Item(31, 3).Value = "Lido DAO"
$ws.Cells.Item(31, 4).Value = 2.34
$ws.Cells.Item(31, 5).Value = 2051128964
$ws.Cells.Item(31, 6).Value = 47872210
$ws.Cells.Item(31, 7).Value = -2.77998

$ws.Cells.Item(32, 2).Value = "TON"
$ws.Cells.Item(32, 3).Value = "Toncoin"
$ws.Cells.Item(32, 4).Value = 1.36
$ws.Cells.Item(32, 5).Value = 1996434252
$ws.Cells.Item(32, 6).Value = 5958813
$ws.Cells.Item(32, 7).Value = 0.73614

$ws.Cells.Item(33, 2).Value = "FIL"
$ws.Cells.Item(33, 3).Value = "Filecoin"
$ws.Cells.Item(33, 4).Value = 4.44
$ws.Cells.Item(33, 5).Value = 1939263133
$ws.Cells.Item(33, 6).Value = 88561125
$ws.Cells.Item(33, 7).Value = -0.58205

$ws.Cells.Item(34, 2).Value = "ICP"
$ws.Cells.Item(34, 3).Value = "Internet Computer"
$ws.Cells.Item(34, 4).Value = 4.16
$ws.Cells.Item(34, 5).Value = 1826130293
$ws.Cells.Item(34, 6).Value = 16158656
$ws.Cells.Item(34, 7).Value = -1.63878

$ws.Cells.Item(35, 2).Value = "ARB"
$ws.Cells.Item(35, 3).Value = "Arbitrum"
$ws.Cells.Item(35, 4).Value = 1.3
$ws.Cells.Item(35, 5).Value = 1662611797
$ws.Cells.Item(35, 6).Value = 287584018
$ws.Cells.Item(35, 7).Value = -0.66048

$ws.Cells.Item(36, 2).Value = "HBAR"
$ws.Cells.Item(36, 3).Value = "Hedera"
$ws.Cells.Item(36, 4).Value = 0.051243
$ws.Cells.Item(36, 5).Value = 1654913060
$ws.Cells.Item(36, 6).Value = 20803854
$ws.Cells.Item(36, 7).Value = -2.64609

$ws.Cells.Item(37, 2).Value = "APT"
$ws.Cells.Item(37, 3).Value = "Aptos"
$ws.Cells.Item(37, 4).Value = 7.58
$ws.Cells.Item(37, 5).Value = 1627826672
$ws.Cells.Item(37, 6).Value = 72234937
$ws.Cells.Item(37, 7).Value = 0.87093

$ws.Cells.Item(38, 2).Value = "CRO"
$ws.Cells.Item(38, 3).Value = "Cronos"
$ws.Cells.Item(38, 4).Value = 0.060809
$ws.Cells.Item(38, 5).Value = 1590196384
$ws.Cells.Item(38, 6).Value = 9419116
$ws.Cells.Item(38, 7).Value = 0.16994

$ws.Cells.Item(39, 2).Value = "QNT"
$ws.Cells.Item(39, 3).Value = "Quant"
$ws.Cells.Item(39, 4).Value = 102.03
$ws.Cells.Item(39, 5).Value = 1483216203
$ws.Cells.Item(39, 6).Value = 12586638
$ws.Cells.Item(39, 7).Value = -0.47577

$ws.Cells.Item(40, 2).Value = "VET"
$ws.Cells.Item(40, 3).Value = "VeChain"
$ws.Cells.Item(40, 4).Value = 0.01973754
$ws.Cells.Item(40, 5).Value = 1436239945
$ws.Cells.Item(40, 6).Value = 37093969
$ws.Cells.Item(40, 7).Value = -1.04918

$ws.Cells.Item(41, 2).Value = "NEAR"
$ws.Cells.Item(41, 3).Value = "NEAR Protocol"
$ws.Cells.Item(41, 4).Value = 1.49
$ws.Cells.Item(41, 5).Value = 1393918332
$ws.Cells.Item(41, 6).Value = 43521535
$ws.Cells.Item(41, 7).Value = -0.82492

$ws.Cells.Item(42, 2).Value = "AAVE"
$ws.Cells.Item(42, 3).Value = "Aave"
$ws.Cells.Item(42, 4).Value = 78.44
$ws.Cells.Item(42, 5).Value = 1136417724
$ws.Cells.Item(42, 6).Value = 84082017
$ws.Cells.Item(42, 7).Value = -4.38826

$ws.Cells.Item(43, 2).Value = "GRT"
$ws.Cells.Item(43, 3).Value = "The Graph"
$ws.Cells.Item(43, 4).Value = 0.121644
$ws.Cells.Item(43, 5).Value = 1104398180
$ws.Cells.Item(43, 6).Value = 27493817
$ws.Cells.Item(43, 7).Value = 0.399

$ws.Cells.Item(44, 2).Value = "FRAX"
$ws.Cells.Item(44, 3).Value = "Frax"
$ws.Cells.Item(44, 4).Value = 1.001
$ws.Cells.Item(44, 5).Value = 1005435036
$ws.Cells.Item(44, 6).Value = 3507306
$ws.Cells.Item(44, 7).Value = 0.29938

$ws.Cells.Item(45, 2).Value = "RETH"
$ws.Cells.Item(45, 3).Value = "Rocket Pool ETH"
$ws.Cells.Item(45, 4).Value = 2085.7
$ws.Cells.Item(45, 5).Value = 998975997
$ws.Cells.Item(45, 6).Value = 1230429
$ws.Cells.Item(45, 7).Value = -0.18914

$ws.Cells.Item(46, 2).Value = "EGLD"
$ws.Cells.Item(46, 3).Value = "MultiversX"
$ws.Cells.Item(46, 4).Value = 37.49
$ws.Cells.Item(46, 5).Value = 964828911
$ws.Cells.Item(46, 6).Value = 10021752
$ws.Cells.Item(46, 7).Value = 1.41532

$ws.Cells.Item(47, 2).Value = "OP"
$ws.Cells.Item(47, 3).Value = "Optimism"
$ws.Cells.Item(47, 4).Value = 1.48
$ws.Cells.Item(47, 5).Value = 953205503
$ws.Cells.Item(47, 6).Value = 119715687
$ws.Cells.Item(47, 7).Value = -2.00554

$ws.Cells.Item(48, 2).Value = "ALGO"
$ws.Cells.Item(48, 3).Value = "Algorand"
$ws.Cells.Item(48, 4).Value = 0.117502
$ws.Cells.Item(48, 5).Value = 906081022
$ws.Cells.Item(48, 6).Value = 32769066
$ws.Cells.Item(48, 7).Value = 2.72119

$ws.Cells.Item(49, 2).Value = "MKR"
$ws.Cells.Item(49, 3).Value = "Maker"
$ws.Cells.Item(49, 4).Value = 993.4
$ws.Cells.Item(49, 5).Value = 897797164
$ws.Cells.Item(49, 6).Value = 82307323
$ws.Cells.Item(49, 7).Value = 12.62316

$ws.Cells.Item(50, 2).Value = "STX"
$ws.Cells.Item(50, 3).Value = "Stacks"
$ws.Cells.Item(50, 4).Value = 0.643779
$ws.Cells.Item(50, 5).Value = 896572036
$ws.Cells.Item(50, 6).Value = 9894914
$ws.Cells.Item(50, 7).Value = -0.8104

$ws.Cells.Item(51, 2).Value = "SNX"
$ws.Cells.Item(51, 3).Value = "Synthetix Network"
$ws.Cells.Item(51, 4).Value = 2.76
$ws.Cells.Item(51, 5).Value = 885309274
$ws.Cells.Item(51, 6).Value = 109913917
$ws.Cells.Item(51, 7).Value = 1.44065

